# edit.ps1 - reproduce the target commit against before.pptx
#
# The commit makes two logical changes to the deck:
#
#   1. The table on slide 6 (the "SOURCES OF FINANCE" slide) gets a new
#      table style (tableStyleId {16815A15-071F-4391-B38B-2CB52E5E059B}
#      -> {9910307B-8CA1-41B6-B9FD-A21C34D8FAEF}).
#
#   2. The presentation's design theme is switched from the custom
#      "Integral" colour scheme back to the built-in "Office Theme"
#      colour scheme (12 theme colours: dk1/lt1/dk2/lt2/accent1-6/
#      hlink/folHlink).

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Table style on slide 6
# -----------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{9910307B-8CA1-41B6-B9FD-A21C34D8FAEF}")

# -----------------------------------------------------------------
# 2) Switch the deck's theme colours to the standard "Office Theme"
#    palette (was the green "Integral" palette).
# -----------------------------------------------------------------
function Set-DeckThemeColor {
    param($colorScheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme

Set-DeckThemeColor $themeColors 1  "000000"   # dk1
Set-DeckThemeColor $themeColors 2  "FFFFFF"   # lt1
Set-DeckThemeColor $themeColors 3  "44546A"   # dk2
Set-DeckThemeColor $themeColors 4  "E7E6E6"   # lt2
Set-DeckThemeColor $themeColors 5  "5B9BD5"   # accent1
Set-DeckThemeColor $themeColors 6  "ED7D31"   # accent2
Set-DeckThemeColor $themeColors 7  "A5A5A5"   # accent3
Set-DeckThemeColor $themeColors 8  "FFC000"   # accent4
Set-DeckThemeColor $themeColors 9  "4472C4"   # accent5
Set-DeckThemeColor $themeColors 10 "70AD47"   # accent6
Set-DeckThemeColor $themeColors 11 "0563C1"   # hlink
Set-DeckThemeColor $themeColors 12 "954F72"   # folHlink
